$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '65.580.24'
$ws.Range('E2').Value = '  +1.64%  '

$ws.Range('D3').Value = '3.400.75'
$ws.Range('E3').Value = '  +1.21%  '

$c = $ws.Range('D4')
$s = $c.Style
$c.Value = '''0.999'
$c.Style = $s
$ws.Range('E4').Value = '  -0.22%  '

$c = $ws.Range('D5')
$s = $c.Style
$c.Value = '''560.56'
$c.Style = $s
$ws.Range('E5').Value = '  -0.37%  '

$c = $ws.Range('D6')
$s = $c.Style
$c.Value = '''176.24'
$c.Style = $s
$ws.Range('E6').Value = '  +0.30%  '

$c = $ws.Range('D7')
$s = $c.Style
$c.Value = '''0.633'
$c.Style = $s
$ws.Range('E7').Value = '  +0.44%  '

$ws.Range('D8').Value = '3.391.95'
$ws.Range('E8').Value = '  +1.26%  '

$c = $ws.Range('D9')
$s = $c.Style
$c.Value = '''0.999'
$c.Style = $s
$ws.Range('E9').Value = '  -0.22%  '

$ws.Range('E10').Value = '  +5.02%  '

$c = $ws.Range('D11')
$s = $c.Style
$c.Value = '''0.641'
$c.Style = $s
$ws.Range('E11').Value = '  +1.38%  '

$c = $ws.Range('D12')
$s = $c.Style
$c.Value = '''53.53'
$c.Style = $s
$ws.Range('E12').Value = '  -2.96%  '

$ws.Range('E13').Value = '  +1.15%  '

$c = $ws.Range('D14')
$s = $c.Style
$c.Value = '''9.23'
$c.Style = $s
$ws.Range('E14').Value = '  +1.46%  '

$ws.Range('D15').Value = '3.932.28'
$ws.Range('E15').Value = '  +0.71%  '

$c = $ws.Range('D16')
$s = $c.Style
$c.Value = '''18.35'
$c.Style = $s
$ws.Range('E16').Value = '  +0.39%  '

$ws.Range('E17').Value = '  +1.48%  '

$ws.Range('D18').Value = '3.402.75'
$ws.Range('E18').Value = '  +1.16%  '

$ws.Range('D19').Value = '65.560.62'
$ws.Range('E19').Value = '  +1.69%  '

$c = $ws.Range('D20')
$s = $c.Style
$c.Value = '''11.86'
$c.Style = $s
$ws.Range('E20').Value = '  +0.16%  '

$c = $ws.Range('D21')
$s = $c.Style
$c.Value = '''1.01'
$c.Style = $s
$ws.Range('E21').Value = '  +1.67%  '

$c = $ws.Range('D22')
$s = $c.Style
$c.Value = '''490.22'
$c.Style = $s
$ws.Range('E22').Value = '  +4.82%  '

$c = $ws.Range('D23')
$s = $c.Style
$c.Value = '''4.95'
$c.Style = $s
$ws.Range('E23').Value = '  -0.82%  '

$c = $ws.Range('D24')
$s = $c.Style
$c.Value = '''4.13'
$c.Style = $s
$ws.Range('E24').Value = '  +0.04%  '

$c = $ws.Range('D25')
$s = $c.Style
$c.Value = '''89.13'
$c.Style = $s
$ws.Range('E25').Value = '  +2.78%  '

$c = $ws.Range('D26')
$s = $c.Style
$c.Value = '''14.20'
$c.Style = $s
$ws.Range('E26').Value = '  +4.89%  '

$c = $ws.Range('D27')
$s = $c.Style
$c.Value = '''2.92'
$c.Style = $s
$ws.Range('E27').Value = '  +2.89%  '

$c = $ws.Range('D28')
$s = $c.Style
$c.Value = '''10.74'
$c.Style = $s
$ws.Range('E28').Value = '  -0.89%  '

$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D29')
$s = $c.Style
$c.Value = '''8.74'
$c.Style = $s
$ws.Range('E29').Value = '  -0.76%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D30')
$s = $c.Style
$c.Value = '''31.47'
$c.Style = $s
$ws.Range('E30').Value = '  +4.31%  '

$c = $ws.Range('D31')
$s = $c.Style
$c.Value = '''6.57'
$c.Style = $s
$ws.Range('E31').Value = '  -1.19%  '

$c = $ws.Range('D32')
$s = $c.Style
$c.Value = '''11.51'
$c.Style = $s
$ws.Range('E32').Value = '  +0.15%  '

$c = $ws.Range('D33')
$s = $c.Style
$c.Value = '''62.71'
$c.Style = $s
$ws.Range('E33').Value = '  +5.41%  '

$c = $ws.Range('D34')
$s = $c.Style
$c.Value = '''576.44'
$c.Style = $s
$ws.Range('E34').Value = '  -0.50%  '

$ws.Range('E35').Value = '  -0.16%  '

$ws.Range('E36').Value = '  +0.03%  '

$ws.Range('E37').Value = '  +5.07%  '

$c = $ws.Range('D38')
$s = $c.Style
$c.Value = '''0.141'
$c.Style = $s
$ws.Range('E38').Value = '  +0.53%  '

$c = $ws.Range('D39')
$s = $c.Style
$c.Value = '''35.99'
$c.Style = $s
$ws.Range('E39').Value = '  +0.11%  '

$ws.Range('E40').Value = '  +0.87%  '

$ws.Range('D41').Value = '0.0₃0744'
$ws.Range('E41').Value = '  -1.47%  '

$ws.Range('D42').Value = '3.119.74'
$ws.Range('E42').Value = '  +0.97%  '

$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D43')
$s = $c.Style
$c.Value = '''2.79'
$c.Style = $s
$ws.Range('E43').Value = '  -1.05%  '

$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D44')
$s = $c.Style
$c.Value = '''0.0418'
$c.Style = $s
$ws.Range('E44').Value = '  +1.38%  '

$c = $ws.Range('D45')
$s = $c.Style
$c.Value = '''0.135'
$c.Style = $s
$ws.Range('E45').Value = '  +1.85%  '

$c = $ws.Range('D46')
$s = $c.Style
$c.Value = '''3.19'
$c.Style = $s
$ws.Range('E46').Value = '  -0.37%  '

$ws.Range('E47').Value = '  -3.28%  '

$c = $ws.Range('D48')
$s = $c.Style
$c.Value = '''0.998'
$c.Style = $s
$ws.Range('E48').Value = '  -0.18%  '

$c = $ws.Range('D49')
$s = $c.Style
$c.Value = '''140.35'
$c.Style = $s
$ws.Range('E49').Value = '  +2.15%  '

$ws.Range('E50').Value = '  -1.57%  '

$c = $ws.Range('D51')
$s = $c.Style
$c.Value = '''8.46'
$c.Style = $s
$ws.Range('E51').Value = '  +0.87%  '
